# Logic tree input file updated
#
# A new "Possible_Problem" row is inserted at row 8 (pushing the former
# rows 8-9 down to rows 9-10). The new row carries the same
# Possible_Problem / percentage-breakdown pair already used in row 4
# (same shared strings), under the "Does driving the vehicle alleviate
# the problem?" question in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8, shifting old rows 8 and 9
# down to rows 9 and 10.
$ws.Rows.Item(8).Insert()

# Column A keeps repeating the same "Does driving..." question used by
# the rows below it.
$ws.Range("A8").Value = $ws.Range("A9").Value2

# New Possible_Problem / percentage-breakdown pair (same text as row 4).
$ws.Range("B8").Value = "Possible_Problem"
$ws.Range("C8").Value = $ws.Range("C4").Value2

# Match the row height Excel auto-computed for this wrapped, long text
# (same as row 4, which holds identical text).
$ws.Rows.Item(8).RowHeight = 409.6

# Restore the selection/scroll position left behind in the saved file.
$ws.Range("B11").Select()
